$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3683, 3867, 4258, 4586, 4586, 4586, 4635, 4635, 4818, 5232, 5232, 5274, 5274, 5274)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
